$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 78.4078284531868
$ws.Cells.Item(1, 2).Value = 182.564449209051
$ws.Cells.Item(1, 3).Value = 67.3528387524899
$ws.Cells.Item(1, 4).Value = 79.6841525843759
$ws.Cells.Item(1, 5).Value = 28.5875222778821
$ws.Cells.Item(1, 6).Value = 87.8914307280869
$ws.Cells.Item(1, 7).Value = 15.1879239898119
$ws.Cells.Item(1, 8).Value = 11.9976429324586
$ws.Cells.Item(1, 9).Value = 60.7411978117848
$ws.Cells.Item(1, 10).Value = 118.26333390468
$ws.Cells.Item(2, 1).Value = 89.9166619823857
$ws.Cells.Item(2, 2).Value = 30.4632081792053
$ws.Cells.Item(2, 3).Value = 58.6066634667137
$ws.Cells.Item(2, 4).Value = 156.719543019645
$ws.Cells.Item(2, 5).Value = 178.465342232243
$ws.Cells.Item(2, 6).Value = 144.7297707874
$ws.Cells.Item(2, 7).Value = 142.85651079512
$ws.Cells.Item(2, 8).Value = 127.029146592612
$ws.Cells.Item(2, 9).Value = 118.892595227292
$ws.Cells.Item(2, 10).Value = 24.179480515504
$ws.Cells.Item(3, 1).Value = 13.2353526601733
$ws.Cells.Item(3, 2).Value = 16.789862335096
$ws.Cells.Item(3, 3).Value = 64.4125092143251
$ws.Cells.Item(3, 4).Value = 91.0139336674539
$ws.Cells.Item(3, 5).Value = 60.0258050765962
$ws.Cells.Item(3, 6).Value = 190.012591327546
$ws.Cells.Item(3, 7).Value = 165.737399629195
$ws.Cells.Item(3, 8).Value = 136.806406703222
$ws.Cells.Item(3, 9).Value = 70.7165700712784
$ws.Cells.Item(3, 10).Value = 90.2644459578509
$ws.Cells.Item(4, 1).Value = 42.051677052887
$ws.Cells.Item(4, 2).Value = 74.0834727296994
$ws.Cells.Item(4, 3).Value = 96.4409285674062
$ws.Cells.Item(4, 4).Value = 84.5790998472735
$ws.Cells.Item(4, 5).Value = 146.582222239385
$ws.Cells.Item(4, 6).Value = 99.0692046932267
$ws.Cells.Item(4, 7).Value = 171.283539184967
$ws.Cells.Item(4, 8).Value = 70.1933742827705
$ws.Cells.Item(4, 9).Value = 44.858923389045
$ws.Cells.Item(4, 10).Value = 195.312695296161
$ws.Cells.Item(5, 1).Value = 21.3417588832517
$ws.Cells.Item(5, 2).Value = 99.9221572186435
$ws.Cells.Item(5, 3).Value = 123.754329664518
$ws.Cells.Item(5, 4).Value = 22.8757538007925
$ws.Cells.Item(5, 5).Value = 122.170106099066
$ws.Cells.Item(5, 6).Value = 170.510997982002
$ws.Cells.Item(5, 7).Value = 154.685959478228
$ws.Cells.Item(5, 8).Value = 153.219697137
$ws.Cells.Item(5, 9).Value = 46.0396748250535
$ws.Cells.Item(5, 10).Value = 123.379056771928
$ws.Cells.Item(6, 1).Value = 166.750805902644
$ws.Cells.Item(6, 2).Value = 54.0387402540253
$ws.Cells.Item(6, 3).Value = 199.319233093094
$ws.Cells.Item(6, 4).Value = 175.136718701169
$ws.Cells.Item(6, 5).Value = 2.42712525763881
$ws.Cells.Item(6, 6).Value = 61.6179661180908
$ws.Cells.Item(6, 7).Value = 118.151939994726
$ws.Cells.Item(6, 8).Value = 176.338905085036
$ws.Cells.Item(6, 9).Value = 19.6583475077796
$ws.Cells.Item(6, 10).Value = 38.5749309503357
$ws.Cells.Item(7, 1).Value = 122.154031098892
$ws.Cells.Item(7, 2).Value = 78.3815172865901
$ws.Cells.Item(7, 3).Value = 141.28107286118
$ws.Cells.Item(7, 4).Value = 170.476751853934
$ws.Cells.Item(7, 5).Value = 76.2116568517926
$ws.Cells.Item(7, 6).Value = 15.8331892526863
$ws.Cells.Item(7, 7).Value = 134.022279611799
$ws.Cells.Item(7, 8).Value = 174.02756361944
$ws.Cells.Item(7, 9).Value = 10.1373207802592
$ws.Cells.Item(7, 10).Value = 79.3961375390161
$ws.Cells.Item(8, 1).Value = 173.446231602433
$ws.Cells.Item(8, 2).Value = 72.6631365123499
$ws.Cells.Item(8, 3).Value = 82.170223203567
$ws.Cells.Item(8, 4).Value = 123.579899931131
$ws.Cells.Item(8, 5).Value = 2.83772163225232
$ws.Cells.Item(8, 6).Value = 113.31319544153
$ws.Cells.Item(8, 7).Value = 93.0355326705778
$ws.Cells.Item(8, 8).Value = 41.5367554135326
$ws.Cells.Item(8, 9).Value = 168.843827568388
$ws.Cells.Item(8, 10).Value = 89.5148070945939
$ws.Cells.Item(9, 1).Value = 35.3266318493181
$ws.Cells.Item(9, 2).Value = 12.517702492195
$ws.Cells.Item(9, 3).Value = 90.7667318781683
$ws.Cells.Item(9, 4).Value = 147.33751329935
$ws.Cells.Item(9, 5).Value = 123.513640055206
$ws.Cells.Item(9, 6).Value = 188.012936798862
$ws.Cells.Item(9, 7).Value = 74.7642396366057
$ws.Cells.Item(9, 8).Value = 121.304209866237
$ws.Cells.Item(9, 9).Value = 82.1519745896347
$ws.Cells.Item(9, 10).Value = 84.9642561212947
$ws.Cells.Item(10, 1).Value = 180.9172646985
$ws.Cells.Item(10, 2).Value = 194.944634099931
$ws.Cells.Item(10, 3).Value = 50.5350267749909
$ws.Cells.Item(10, 4).Value = 6.28399243870936
$ws.Cells.Item(10, 5).Value = 73.1586641972692
$ws.Cells.Item(10, 6).Value = 142.960241596662
$ws.Cells.Item(10, 7).Value = 158.641084357463
$ws.Cells.Item(10, 8).Value = 153.277577810584
$ws.Cells.Item(10, 9).Value = 146.664096949
$ws.Cells.Item(10, 10).Value = 106.33691684638
$ws.Cells.Item(11, 1).Value = 36.4887183702033
$ws.Cells.Item(11, 2).Value = 180.658395858788
$ws.Cells.Item(11, 3).Value = 143.082376356741
$ws.Cells.Item(11, 4).Value = 166.643537286037
$ws.Cells.Item(11, 5).Value = 149.932825169495
$ws.Cells.Item(11, 6).Value = 94.0876693902946
$ws.Cells.Item(11, 7).Value = 171.868517050458
$ws.Cells.Item(11, 8).Value = 75.7393331619628
$ws.Cells.Item(11, 9).Value = 172.298997068917
$ws.Cells.Item(11, 10).Value = 89.1139298161091
$ws.Cells.Item(12, 1).Value = 168.582433447513
$ws.Cells.Item(12, 2).Value = 76.6151845811937
$ws.Cells.Item(12, 3).Value = 7.49507751664849
$ws.Cells.Item(12, 4).Value = 130.143540413186
$ws.Cells.Item(12, 5).Value = 3.24829910101756
$ws.Cells.Item(12, 6).Value = 109.636328606697
$ws.Cells.Item(12, 7).Value = 187.614785408422
$ws.Cells.Item(12, 8).Value = 193.94355956183
$ws.Cells.Item(12, 9).Value = 46.9631117987275
$ws.Cells.Item(12, 10).Value = 88.1987200529309
$ws.Cells.Item(13, 1).Value = 141.068949616081
$ws.Cells.Item(13, 2).Value = 12.7180697455621
$ws.Cells.Item(13, 3).Value = 91.8755890298056
$ws.Cells.Item(13, 4).Value = 125.173064658965
$ws.Cells.Item(13, 5).Value = 98.4788728405157
$ws.Cells.Item(13, 6).Value = 178.501597502502
$ws.Cells.Item(13, 7).Value = 22.128109737359
$ws.Cells.Item(13, 8).Value = 75.8862307648576
$ws.Cells.Item(13, 9).Value = 50.4212357338617
$ws.Cells.Item(13, 10).Value = 59.8774800355907
$ws.Cells.Item(14, 1).Value = 154.672111084066
$ws.Cells.Item(14, 2).Value = 139.757954859993
$ws.Cells.Item(14, 3).Value = 94.8726584645327
$ws.Cells.Item(14, 4).Value = 62.5069107220075
$ws.Cells.Item(14, 5).Value = 53.0260887243906
$ws.Cells.Item(14, 6).Value = 54.66823599053
$ws.Cells.Item(14, 7).Value = 69.4353261354544
$ws.Cells.Item(14, 8).Value = 124.123194592131
$ws.Cells.Item(14, 9).Value = 197.404688129856
$ws.Cells.Item(14, 10).Value = 29.4259706649119
$ws.Cells.Item(15, 1).Value = 16.1444197484033
$ws.Cells.Item(15, 2).Value = 199.024906474643
$ws.Cells.Item(15, 3).Value = 149.00521279732
$ws.Cells.Item(15, 4).Value = 193.038044773526
$ws.Cells.Item(15, 5).Value = 116.381822673782
$ws.Cells.Item(15, 6).Value = 104.302080117307
$ws.Cells.Item(15, 7).Value = 187.449556583282
$ws.Cells.Item(15, 8).Value = 120.391486361805
$ws.Cells.Item(15, 9).Value = 3.0356933376918
$ws.Cells.Item(15, 10).Value = 163.522335590572
$ws.Cells.Item(16, 1).Value = 155.34545618824
$ws.Cells.Item(16, 2).Value = 164.697524795634
$ws.Cells.Item(16, 3).Value = 106.314466011857
$ws.Cells.Item(16, 4).Value = 58.4653768960691
$ws.Cells.Item(16, 5).Value = 165.267967230299
$ws.Cells.Item(16, 6).Value = 23.7706486246412
$ws.Cells.Item(16, 7).Value = 88.7828068289826
$ws.Cells.Item(16, 8).Value = 17.9093116977761
$ws.Cells.Item(16, 9).Value = 68.1646644455216
$ws.Cells.Item(16, 10).Value = 171.431227666992
$ws.Cells.Item(17, 1).Value = 71.9595596529355
$ws.Cells.Item(17, 2).Value = 95.9822862856008
$ws.Cells.Item(17, 3).Value = 25.3180974281012
$ws.Cells.Item(17, 4).Value = 112.421517033326
$ws.Cells.Item(17, 5).Value = 134.441818732043
$ws.Cells.Item(17, 6).Value = 28.8244785875196
$ws.Cells.Item(17, 7).Value = 181.742526116661
$ws.Cells.Item(17, 8).Value = 144.988166794641
$ws.Cells.Item(17, 9).Value = 77.1174516887951
$ws.Cells.Item(17, 10).Value = 148.580063110488
$ws.Cells.Item(18, 1).Value = 40.2010024712426
$ws.Cells.Item(18, 2).Value = 63.4915908162909
$ws.Cells.Item(18, 3).Value = 196.538871431974
$ws.Cells.Item(18, 4).Value = 17.5371411338156
$ws.Cells.Item(18, 5).Value = 72.0543003045275
$ws.Cells.Item(18, 6).Value = 142.044043141438
$ws.Cells.Item(18, 7).Value = 63.712856948242
$ws.Cells.Item(18, 8).Value = 98.83754425628
$ws.Cells.Item(18, 9).Value = 8.79124198518286
$ws.Cells.Item(18, 10).Value = 194.176792723209
$ws.Cells.Item(19, 1).Value = 191.05204091922
$ws.Cells.Item(19, 2).Value = 101.736623375554
$ws.Cells.Item(19, 3).Value = 72.8505374271658
$ws.Cells.Item(19, 4).Value = 86.8989001432894
$ws.Cells.Item(19, 5).Value = 104.532023847351
$ws.Cells.Item(19, 6).Value = 189.974586288433
$ws.Cells.Item(19, 7).Value = 33.4434888481365
$ws.Cells.Item(19, 8).Value = 36.4072815684636
$ws.Cells.Item(19, 9).Value = 97.238943491708
$ws.Cells.Item(19, 10).Value = 29.2554400997495
$ws.Cells.Item(20, 1).Value = 165.885429161547
$ws.Cells.Item(20, 2).Value = 51.5260144376783
$ws.Cells.Item(20, 3).Value = 55.9585301466093
$ws.Cells.Item(20, 4).Value = 25.9734604628633
$ws.Cells.Item(20, 5).Value = 157.466411011976
$ws.Cells.Item(20, 6).Value = 120.162133462803
$ws.Cells.Item(20, 7).Value = 173.706809046542
$ws.Cells.Item(20, 8).Value = 36.5836957639939
$ws.Cells.Item(20, 9).Value = 58.5962260414829
$ws.Cells.Item(20, 10).Value = 87.5573440862621
